# Apply updated profit/price figures to each job sheet's leve-profit table (H:N columns).
# Values were refreshed by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 40008.5
$ws.Range("I21").Value = 10017
$ws.Range("J21").Value = 70000
$ws.Range("K21").Value = 10017
$ws.Range("L21").Value = 70000
$ws.Range("M21").Value = -9549
$ws.Range("N21").Value = -70936
$ws.Range("H23").Value = 40008.5
$ws.Range("I23").Value = 10017
$ws.Range("J23").Value = 70000
$ws.Range("K23").Value = 10017
$ws.Range("L23").Value = 70000
$ws.Range("M23").Value = -9783
$ws.Range("N23").Value = -70468
$ws.Range("H129").Value = 862.14
$ws.Range("J129").Value = 893.5106
$ws.Range("L129").Value = 2680.5318
$ws.Range("N129").Value = -12680.5318
$ws.Range("H132").Value = 406130.25
$ws.Range("I132").Value = 6397.619
$ws.Range("J132").Value = 2504726.5
$ws.Range("K132").Value = 19192.857
$ws.Range("L132").Value = 7514179.5
$ws.Range("M132").Value = -16662.857
$ws.Range("N132").Value = -7519239.5
$ws.Range("H137").Value = 2166901.5
$ws.Range("I137").Value = 3664304
$ws.Range("J137").Value = 3986.6667
$ws.Range("K137").Value = 10992912
$ws.Range("L137").Value = 11960.0001
$ws.Range("M137").Value = -10990362
$ws.Range("N137").Value = -17060.0001
$ws.Range("H138").Value = 5769.44
$ws.Range("J138").Value = 7543.946
$ws.Range("L138").Value = 22631.838
$ws.Range("N138").Value = -32911.838

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 829.1
$ws.Range("I2").Value = 652.4706
$ws.Range("J2").Value = 1830
$ws.Range("K2").Value = 652.4706
$ws.Range("L2").Value = 1830
$ws.Range("M2").Value = -539.4706
$ws.Range("N2").Value = -2056
$ws.Range("H61").Value = 1615.25
$ws.Range("I61").Value = 1508.75
$ws.Range("J61").Value = 1828.25
$ws.Range("K61").Value = 1508.75
$ws.Range("L61").Value = 1828.25
$ws.Range("M61").Value = -1296.75
$ws.Range("N61").Value = -2252.25
$ws.Range("H74").Value = 4807.654
$ws.Range("I74").Value = 5266.722
$ws.Range("J74").Value = 3774.75
$ws.Range("K74").Value = 5266.722
$ws.Range("L74").Value = 3774.75
$ws.Range("M74").Value = -4392.722
$ws.Range("N74").Value = -5522.75
$ws.Range("H77").Value = 4807.654
$ws.Range("I77").Value = 5266.722
$ws.Range("J77").Value = 3774.75
$ws.Range("K77").Value = 26333.61
$ws.Range("L77").Value = 18873.75
$ws.Range("M77").Value = -21965.61
$ws.Range("N77").Value = -27609.75
$ws.Range("H116").Value = 829.1
$ws.Range("I116").Value = 652.4706
$ws.Range("J116").Value = 1830
$ws.Range("K116").Value = 652.4706
$ws.Range("L116").Value = 1830
$ws.Range("M116").Value = 1641.5294
$ws.Range("N116").Value = -6418
$ws.Range("H122").Value = 2077.7368
$ws.Range("I122").Value = 1105.25
$ws.Range("J122").Value = 3744.8572
$ws.Range("K122").Value = 3315.75
$ws.Range("L122").Value = 11234.5716
$ws.Range("M122").Value = -865.75
$ws.Range("N122").Value = -16134.5716
$ws.Range("H132").Value = 2287.818
$ws.Range("I132").Value = 1222.2667
$ws.Range("J132").Value = 4571.143
$ws.Range("K132").Value = 3666.800099999999
$ws.Range("L132").Value = 13713.429
$ws.Range("M132").Value = -1136.800099999999
$ws.Range("N132").Value = -18773.429
$ws.Range("H136").Value = 1615.25
$ws.Range("I136").Value = 1508.75
$ws.Range("J136").Value = 1828.25
$ws.Range("K136").Value = 4526.25
$ws.Range("L136").Value = 5484.75
$ws.Range("M136").Value = -1976.25
$ws.Range("N136").Value = -10584.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 829.1
$ws.Range("I3").Value = 652.4706
$ws.Range("J3").Value = 1830
$ws.Range("K3").Value = 652.4706
$ws.Range("L3").Value = 1830
$ws.Range("M3").Value = -538.4706
$ws.Range("N3").Value = -2058
$ws.Range("H20").Value = 24961.715
$ws.Range("I20").Value = 8010
$ws.Range("J20").Value = 27787
$ws.Range("K20").Value = 8010
$ws.Range("L20").Value = 27787
$ws.Range("M20").Value = -7763
$ws.Range("N20").Value = -28281
$ws.Range("H107").Value = 2006.9524
$ws.Range("I107").Value = 1906.5883
$ws.Range("J107").Value = 2433.5
$ws.Range("K107").Value = 1906.5883
$ws.Range("L107").Value = 2433.5
$ws.Range("M107").Value = 13.41170000000011
$ws.Range("N107").Value = -6273.5
$ws.Range("H134").Value = 2873.2
$ws.Range("I134").Value = 2108.6667
$ws.Range("J134").Value = 4657.1113
$ws.Range("K134").Value = 6326.000100000001
$ws.Range("L134").Value = 13971.3339
$ws.Range("M134").Value = -3791.000100000001
$ws.Range("N134").Value = -19041.3339

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2494.9678
$ws.Range("I31").Value = 1101.9131
$ws.Range("J31").Value = 6500
$ws.Range("K31").Value = 1101.9131
$ws.Range("L31").Value = 6500
$ws.Range("M31").Value = -806.9131
$ws.Range("N31").Value = -7090
$ws.Range("H34").Value = 2494.9678
$ws.Range("I34").Value = 1101.9131
$ws.Range("J34").Value = 6500
$ws.Range("K34").Value = 1101.9131
$ws.Range("L34").Value = 6500
$ws.Range("M34").Value = -899.9131
$ws.Range("N34").Value = -6904
$ws.Range("H58").Value = 3079.322
$ws.Range("I58").Value = 1862.0834
$ws.Range("J58").Value = 8390.909
$ws.Range("K58").Value = 1862.0834
$ws.Range("L58").Value = 8390.909
$ws.Range("M58").Value = -1659.0834
$ws.Range("N58").Value = -8796.909
$ws.Range("H94").Value = 929.3200000000001
$ws.Range("J94").Value = 1265.8462
$ws.Range("L94").Value = 1265.8462
$ws.Range("N94").Value = -2167.8462
$ws.Range("H132").Value = 5637.4116
$ws.Range("I132").Value = 6140.5
$ws.Range("J132").Value = 5190.222
$ws.Range("K132").Value = 18421.5
$ws.Range("L132").Value = 15570.666
$ws.Range("M132").Value = -15891.5
$ws.Range("N132").Value = -20630.666
$ws.Range("H134").Value = 2556.923
$ws.Range("I134").Value = 1224
$ws.Range("J134").Value = 7000
$ws.Range("K134").Value = 3672
$ws.Range("L134").Value = 21000
$ws.Range("M134").Value = -1137
$ws.Range("N134").Value = -26070
$ws.Range("H136").Value = 3079.322
$ws.Range("I136").Value = 1862.0834
$ws.Range("J136").Value = 8390.909
$ws.Range("K136").Value = 5586.2502
$ws.Range("L136").Value = 25172.727
$ws.Range("M136").Value = -3036.2502
$ws.Range("N136").Value = -30272.727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3907072
$ws.Range("I113").Value = 690.64703
$ws.Range("J113").Value = 8334304
$ws.Range("K113").Value = 2071.94109
$ws.Range("L113").Value = 25002912
$ws.Range("M113").Value = 98.0589100000002
$ws.Range("N113").Value = -25007252
$ws.Range("H122").Value = 2350.3962
$ws.Range("I122").Value = 864.9286
$ws.Range("J122").Value = 2883.641
$ws.Range("K122").Value = 7784.3574
$ws.Range("L122").Value = 25952.769
$ws.Range("M122").Value = -5334.3574
$ws.Range("N122").Value = -30852.769

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5581.8335
$ws.Range("I70").Value = 5260.485
$ws.Range("K70").Value = 5260.485
$ws.Range("M70").Value = -4990.485
$ws.Range("H73").Value = 5581.8335
$ws.Range("I73").Value = 5260.485
$ws.Range("K73").Value = 5260.485
$ws.Range("M73").Value = -4324.485
$ws.Range("H97").Value = 1456.4286
$ws.Range("I97").Value = 1450
$ws.Range("K97").Value = 1450
$ws.Range("M97").Value = -954
$ws.Range("H132").Value = 5142.5713
$ws.Range("I132").Value = 2250
$ws.Range("J132").Value = 8999.333000000001
$ws.Range("K132").Value = 6750
$ws.Range("L132").Value = 26997.999
$ws.Range("M132").Value = -4220
$ws.Range("N132").Value = -32057.999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5124.1333
$ws.Range("I7").Value = 4021.8333
$ws.Range("J7").Value = 5859
$ws.Range("K7").Value = 4021.8333
$ws.Range("L7").Value = 5859
$ws.Range("M7").Value = -3909.8333
$ws.Range("N7").Value = -6083
$ws.Range("H40").Value = 5970.905
$ws.Range("I40").Value = 4531.9375
$ws.Range("J40").Value = 10575.6
$ws.Range("K40").Value = 4531.9375
$ws.Range("L40").Value = 10575.6
$ws.Range("M40").Value = -4395.9375
$ws.Range("N40").Value = -10847.6
$ws.Range("H122").Value = 5152.95
$ws.Range("I122").Value = 2325.5
$ws.Range("K122").Value = 6976.5
$ws.Range("M122").Value = -4526.5
$ws.Range("H126").Value = 5124.1333
$ws.Range("I126").Value = 4021.8333
$ws.Range("J126").Value = 5859
$ws.Range("K126").Value = 12065.4999
$ws.Range("L126").Value = 17577
$ws.Range("M126").Value = -9595.499899999999
$ws.Range("N126").Value = -22517
$ws.Range("H132").Value = 7859.75
$ws.Range("I132").Value = 3207.7144
$ws.Range("J132").Value = 11478
$ws.Range("K132").Value = 9623.143199999999
$ws.Range("L132").Value = 34434
$ws.Range("M132").Value = -7093.143199999999
$ws.Range("N132").Value = -39494
$ws.Range("H136").Value = 6138.2354
$ws.Range("I136").Value = 2430
$ws.Range("J136").Value = 7683.3335
$ws.Range("K136").Value = 7290
$ws.Range("L136").Value = 23050.0005
$ws.Range("M136").Value = -4740
$ws.Range("N136").Value = -28150.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4387.533
$ws.Range("I122").Value = 2645.3333
$ws.Range("J122").Value = 7000.8335
$ws.Range("K122").Value = 7935.999899999999
$ws.Range("L122").Value = 21002.5005
$ws.Range("M122").Value = -5485.999899999999
$ws.Range("N122").Value = -25902.5005
$ws.Range("H132").Value = 30307190
$ws.Range("I132").Value = 2800
$ws.Range("J132").Value = 41671336
$ws.Range("K132").Value = 8400
$ws.Range("L132").Value = 125014008
$ws.Range("M132").Value = -5870
$ws.Range("N132").Value = -125019068
$ws.Range("H136").Value = 4131.0713
$ws.Range("I136").Value = 2103.0435
$ws.Range("J136").Value = 13460
$ws.Range("K136").Value = 6309.130500000001
$ws.Range("L136").Value = 40380
$ws.Range("M136").Value = -3759.130500000001
$ws.Range("N136").Value = -45480
